$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was inserted at row 83 (pushing every
# subsequent record down by one row, through the former row 165 which is
# now row 166). Insert a blank row at 83 first so everything below shifts
# down intact, then clone the row that lands at 84 (the former row 83)
# back into the freshly inserted row 83, and finally overwrite the two
# cells that actually hold new data for this new record (Fecha / Volumen).
$ws.Rows("83:83").Insert()

$ws.Range("A84:R84").Copy($ws.Range("A83:R83"))

$ws.Range("D83").Value = 44810
$ws.Range("J83").Value = 3000
